# Update Mark Du's timecard entry from "19h 37m" to "22h 37m"
# (the commit message "Updating website and timecard" / diff shows the
# B6 cell value for Mark Du's row changing on Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B6").Value = "22h 37m"

# Author's last active cell/selection moved to N10 before saving.
$ws.Range("N10").Select()
